$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay formatted/stored as plain text,
# since many look like numbers (e.g. "1.00", "601.11") and Excel
# would otherwise coerce them into numeric values, losing formatting
# (trailing zeros) and introducing floating point artifacts.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.190.66"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.612.75"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.11"
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.09"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.610.63"
$ws.Range("E7").Value = "  +2.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.620"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.202"
$ws.Range("E10").Value = "  +7.31%  "
$ws.Range("E11").Value = "  +7.68%  "
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.41"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.187.57"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.51"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "621.20"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.605.49"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.259.73"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.59"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.895"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.35"
$ws.Range("E23").Value = "  -16.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.20"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.11"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.81"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.39"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.38"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.58"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.11"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.31"
$ws.Range("E33").Value = "  +5.47%  "
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "627.86"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.76"
$ws.Range("E36").Value = "  +8.27%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.94"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  +6.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.70"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.143"
$ws.Range("E42").Value = "  +6.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.412.25"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0726"
$ws.Range("E45").Value = "  +3.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.32"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.87"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  -0.04%  "

# Row 46/47 swap (Fetch.AI <-> ThetaToken)
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.01"
$ws.Range("E46").Value = "  +9.88%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  +6.92%  "
